# Added image padding and teeth locations adjustment:
# Insert a new header row above the data, labeling each of the 14 data
# columns with its tooth location ("1-7" .. "1-1", "2-1" .. "2-7"),
# carrying over the right-hand border already used by the table and
# marking the header cells as quoted text (so "1-7" etc. stay text and
# are not reinterpreted as dates).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing data down one row to make room for the new header row.
$ws.Rows("1:1").Insert()

# Grab the "right border only" cell format (already used on column N of the
# data rows) and stamp it across the new header row so the inserted cells
# reuse the existing border style instead of inventing a new one.
$ws.Range("N3").Copy() | Out-Null
$ws.Range("A1:N1").PasteSpecial(-4122) | Out-Null

# Column headers, left-to-right across A1:N1.
$headers = @("1-7", "1-6", "1-5", "1-4", "1-3", "1-2", "1-1", "2-1", "2-2", "2-3", "2-4", "2-5", "2-6", "2-7")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    # Leading apostrophe forces text (quote-prefixed) entry, matching labels
    # like "1-7" that would otherwise be parsed as a date.
    $cell.Value = "'" + $headers[$i]
}

$ws.Range("B18").Select() | Out-Null
